# Apply edits to "Test Cases" sheet (B suite.xlsx) per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- E2: Results column SKIP -> FAIL ---
$ws.Range("E2").Value = "FAIL"

# --- Remove now-unused placeholder rows 90:99 (only had stray C-column styling) ---
$ws.Range("A90:E99").EntireRow.Delete()

# Reference row 84 holds the canonical formatting (border + fill-less + wrap where needed)
# that rows 85-89 should adopt; copy per-column formats from it, then set final values.

# --- Row 85 ---
$ws.Range("A84").Copy()
$ws.Range("A85").PasteSpecial(-4122)
$ws.Range("A85").Value = "TestCase_B84"
$ws.Range("B84").Copy()
$ws.Range("B85").PasteSpecial(-4122)
$ws.Range("B85").Value = "OPQA-613"
$ws.Range("C84").Copy()
$ws.Range("C85").PasteSpecial(-4122)
$ws.Range("C85").Value = "Verify that following fields get displayed correctly for a patent in ALL search results page:`na)Title`nb)Inventors`nc)Assignees`nd)Patent number`ne)Publication date`nf)Times cited count`ng)Comments count"
$ws.Range("D84").Copy()
$ws.Range("D85").PasteSpecial(-4122)
$ws.Range("D85").Value = "Y"
$ws.Range("E84").Copy()
$ws.Range("E85").PasteSpecial(-4122)
$ws.Range("E85").ClearContents()
$ws.Range("A85").RowHeight = 120

# --- Row 86 ---
$ws.Range("A84").Copy()
$ws.Range("A86").PasteSpecial(-4122)
$ws.Range("A86").Value = "TestCase_B85"
$ws.Range("B84").Copy()
$ws.Range("B86").PasteSpecial(-4122)
$ws.Range("B86").Value = "OPQA-614"
$ws.Range("C84").Copy()
$ws.Range("C86").PasteSpecial(-4122)
$ws.Range("C86").Value = "Verify that following fields get displayed correctly for a patent in PATENTS search results page:`na)Title`nb)Inventors`nc)Assignees`nd)Patent number`ne)Publication date`nf)Times cited count`ng)Comments count"
$ws.Range("D84").Copy()
$ws.Range("D86").PasteSpecial(-4122)
$ws.Range("D86").Value = "Y"
$ws.Range("E84").Copy()
$ws.Range("E86").PasteSpecial(-4122)
$ws.Range("E86").ClearContents()
$ws.Range("A86").RowHeight = 120

# --- Row 87 ---
$ws.Range("A84").Copy()
$ws.Range("A87").PasteSpecial(-4122)
$ws.Range("A87").Value = "TestCase_B86"
$ws.Range("B84").Copy()
$ws.Range("B87").PasteSpecial(-4122)
$ws.Range("B87").Value = "OPQA-562"
$ws.Range("C84").Copy()
$ws.Range("C87").PasteSpecial(-4122)
$ws.Range("C87").Value = "Verify that following fields get displayed correctly for an patent in record view page:`na)Title`nb)Inventors`nc)Assignees`nd)Publication Date`ne)Publication Number`nf)Times Cited count`ng)Cited patents count`nh)Cited Articles count`ng)Comments count`ni)Abstract`nj)IPC Codes`nk)DETAILS link"
$ws.Range("D84").Copy()
$ws.Range("D87").PasteSpecial(-4122)
$ws.Range("D87").Value = "Y"
$ws.Range("E84").Copy()
$ws.Range("E87").PasteSpecial(-4122)
$ws.Range("E87").ClearContents()
$ws.Range("A87").RowHeight = 195

# --- Row 88 ---
$ws.Range("A84").Copy()
$ws.Range("A88").PasteSpecial(-4122)
$ws.Range("A88").Value = "TestCase_B87"
$ws.Range("B84").Copy()
$ws.Range("B88").PasteSpecial(-4122)
$ws.Range("B88").Value = "OPQA-567"
$ws.Range("C84").Copy()
$ws.Range("C88").PasteSpecial(-4122)
$ws.Range("C88").Value = "Verify that record view page of a patent gets displayed when user clicks on article title in ALL search results page"
$ws.Range("D84").Copy()
$ws.Range("D88").PasteSpecial(-4122)
$ws.Range("D88").Value = "Y"
$ws.Range("E84").Copy()
$ws.Range("E88").PasteSpecial(-4122)
$ws.Range("E88").ClearContents()

# --- Row 89 ---
$ws.Range("A84").Copy()
$ws.Range("A89").PasteSpecial(-4122)
$ws.Range("A89").Value = "TestCase_B88"
$ws.Range("B84").Copy()
$ws.Range("B89").PasteSpecial(-4122)
$ws.Range("B89").Value = "OPQA-573"
$ws.Range("C84").Copy()
$ws.Range("C89").PasteSpecial(-4122)
$ws.Range("C89").Value = "Verify that record view page of a patent gets displayed when user clicks a patent title in PATENTS search results page"
$ws.Range("D84").Copy()
$ws.Range("D89").PasteSpecial(-4122)
$ws.Range("D89").Value = "Y"
$ws.Range("E84").Copy()
$ws.Range("E89").PasteSpecial(-4122)
$ws.Range("E89").ClearContents()

$excel.CutCopyMode = 0

Write-Host "Edits applied"